$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 135
$ws.Range("H135").Value = 1523.85
$ws.Range("J135").Value = 5647.5
$ws.Range("L135").Value = 50827.5
$ws.Range("N135").Value = -55897.5

# Row 137
$ws.Range("H137").Value = 6675.8096
$ws.Range("I137").Value = 12598.223
$ws.Range("J137").Value = 2234
$ws.Range("K137").Value = 37794.669
$ws.Range("L137").Value = 6702
$ws.Range("M137").Value = -35244.669
$ws.Range("N137").Value = -11802

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2102881
$ws.Range("I2").Value = 2297.5833
$ws.Range("J2").Value = 14706382
$ws.Range("K2").Value = 2297.5833
$ws.Range("L2").Value = 14706382
$ws.Range("M2").Value = -2184.5833
$ws.Range("N2").Value = -14706608

# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("N13").Value = 0

# Row 22
$ws.Range("H22").Value = 8960
$ws.Range("I22").Value = 925
$ws.Range("J22").Value = 16995
$ws.Range("K22").Value = 925
$ws.Range("L22").Value = 16995
$ws.Range("M22").Value = -626
$ws.Range("N22").Value = -17593

# Row 32
$ws.Range("H32").Value = 10833.216
$ws.Range("I32").Value = 6278.3335
$ws.Range("J32").Value = 25636.584
$ws.Range("K32").Value = 6278.3335
$ws.Range("L32").Value = 25636.584
$ws.Range("M32").Value = -5991.3335
$ws.Range("N32").Value = -26210.584

# Row 61
$ws.Range("H61").Value = 2556.913
$ws.Range("I61").Value = 2012.4193
$ws.Range("J61").Value = 3682.2
$ws.Range("K61").Value = 2012.4193
$ws.Range("L61").Value = 3682.2
$ws.Range("M61").Value = -1800.4193
$ws.Range("N61").Value = -4106.2

# Row 74
$ws.Range("H74").Value = 4551226
$ws.Range("I74").Value = 7143548
$ws.Range("J74").Value = 14662.375
$ws.Range("K74").Value = 7143548
$ws.Range("L74").Value = 14662.375
$ws.Range("M74").Value = -7142674
$ws.Range("N74").Value = -16410.375

# Row 77
$ws.Range("H77").Value = 4551226
$ws.Range("I77").Value = 7143548
$ws.Range("J77").Value = 14662.375
$ws.Range("K77").Value = 35717740
$ws.Range("L77").Value = 73311.875
$ws.Range("M77").Value = -35713372
$ws.Range("N77").Value = -82047.875

# Row 116
$ws.Range("H116").Value = 2102881
$ws.Range("I116").Value = 2297.5833
$ws.Range("J116").Value = 14706382
$ws.Range("K116").Value = 2297.5833
$ws.Range("L116").Value = 14706382
$ws.Range("M116").Value = -3.583299999999781
$ws.Range("N116").Value = -14710970

# Row 136
$ws.Range("H136").Value = 2556.913
$ws.Range("I136").Value = 2012.4193
$ws.Range("J136").Value = 3682.2
$ws.Range("K136").Value = 6037.257900000001
$ws.Range("L136").Value = 11046.6
$ws.Range("M136").Value = -3487.257900000001
$ws.Range("N136").Value = -16146.6

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2102881
$ws.Range("I3").Value = 2297.5833
$ws.Range("J3").Value = 14706382
$ws.Range("K3").Value = 2297.5833
$ws.Range("L3").Value = 14706382
$ws.Range("M3").Value = -2183.5833
$ws.Range("N3").Value = -14706610

# Row 63
$ws.Range("H63").Value = 56060
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 56060
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").Value = 56060
$ws.Range("N63").Value = -57432

# Row 66
$ws.Range("H66").Value = 56060
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 56060
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 168180
$ws.Range("M66").Value = 0
$ws.Range("N66").Value = -175044

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1838.742
$ws.Range("I31").Value = 1728.875
$ws.Range("J31").Value = 1876.9565
$ws.Range("K31").Value = 1728.875
$ws.Range("L31").Value = 1876.9565
$ws.Range("M31").Value = -1433.875
$ws.Range("N31").Value = -2466.9565

# Row 34
$ws.Range("H34").Value = 1838.742
$ws.Range("I34").Value = 1728.875
$ws.Range("J34").Value = 1876.9565
$ws.Range("K34").Value = 1728.875
$ws.Range("L34").Value = 1876.9565
$ws.Range("M34").Value = -1526.875
$ws.Range("N34").Value = -2280.9565

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0

# Row 132
$ws.Range("H132").Value = 3655.72
$ws.Range("I132").Value = 3899.7334
$ws.Range("J132").Value = 3289.7
$ws.Range("K132").Value = 11699.2002
$ws.Range("L132").Value = 9869.099999999999
$ws.Range("M132").Value = -9169.200199999999
$ws.Range("N132").Value = -14929.1

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 541.17645
$ws.Range("I5").Value = 516.6667
$ws.Range("K5").Value = 1550.0001
$ws.Range("M5").Value = -1438.0001

# Row 21
$ws.Range("H21").Value = 1350
$ws.Range("I21").Value = 200
$ws.Range("J21").Value = 1514.2858
$ws.Range("K21").Value = 600
$ws.Range("L21").Value = 4542.857400000001
$ws.Range("M21").Value = -427
$ws.Range("N21").Value = -4888.857400000001

# Row 135
$ws.Range("H135").Value = 541.17645
$ws.Range("I135").Value = 516.6667
$ws.Range("K135").Value = 4650.0003
$ws.Range("M135").Value = -2115.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 66.2
$ws.Range("I2").Value = 43.666668
$ws.Range("K2").Value = 43.666668
$ws.Range("M2").Value = 69.333332

# Row 132
$ws.Range("H132").Value = 3764.1667
$ws.Range("I132").Value = 3707.25
$ws.Range("J132").Value = 3878
$ws.Range("K132").Value = 11121.75
$ws.Range("L132").Value = 11634
$ws.Range("M132").Value = -8591.75
$ws.Range("N132").Value = -16694

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 6257500
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 6257500
$ws.Range("K2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("M2").Value = 6257500
$ws.Range("N2").Value = -6257724

# Row 40
$ws.Range("H40").Value = 4172.727
$ws.Range("I40").Value = 4367.8
$ws.Range("J40").Value = 2222
$ws.Range("K40").Value = 4367.8
$ws.Range("L40").Value = 2222
$ws.Range("M40").Value = -4231.8
$ws.Range("N40").Value = -2494

# Row 136
$ws.Range("H136").Value = 3015.8823
$ws.Range("I136").Value = 1882.909
$ws.Range("K136").Value = 5648.727000000001
$ws.Range("M136").Value = -3098.727000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Range("H24").Value = 4692.8335
$ws.Range("J24").Value = 4692.8335
$ws.Range("L24").Value = 4692.8335
$ws.Range("N24").Value = -5152.8335

# Row 54
$ws.Range("H54").Value = 11291.333
$ws.Range("I54").Value = 9500
$ws.Range("J54").Value = 11454.182
$ws.Range("K54").Value = 9500
$ws.Range("L54").Value = 11454.182
$ws.Range("M54").Value = -8980
$ws.Range("N54").Value = -12494.182

# Row 136
$ws.Range("H136").Value = 27144.11
$ws.Range("I136").Value = 9069.25
$ws.Range("J136").Value = 46862.137
$ws.Range("K136").Value = 27207.75
$ws.Range("L136").Value = 140586.411
$ws.Range("M136").Value = -24657.75
$ws.Range("N136").Value = -145686.411
